$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H40").Value = 11445.454
$ws_ALC.Range("I40").Value = 21280.4
$ws_ALC.Range("J40").Value = 3249.6667
$ws_ALC.Range("K40").Value = 21280.4
$ws_ALC.Range("L40").Value = 3249.6667
$ws_ALC.Range("M40").Value = -21105.4
$ws_ALC.Range("N40").Value = -3599.6667

$ws_ALC.Range("H53").Value = 84
$ws_ALC.Range("I53").Value = 55.916668
$ws_ALC.Range("K53").Value = 55.916668
$ws_ALC.Range("M53").Value = 581.083332

$ws_ALC.Range("H76").Value = 3006.4333
$ws_ALC.Range("J76").Value = 3040
$ws_ALC.Range("L76").Value = 3040
$ws_ALC.Range("N76").Value = -3670

$ws_ALC.Range("H79").Value = 3006.4333
$ws_ALC.Range("J79").Value = 3040
$ws_ALC.Range("L79").Value = 3040
$ws_ALC.Range("N79").Value = -5224

$ws_ALC.Range("H80").Value = 60607056
$ws_ALC.Range("I80").Value = 66667640
$ws_ALC.Range("J80").Value = 1200
$ws_ALC.Range("K80").Value = 200002920
$ws_ALC.Range("L80").Value = 3600
$ws_ALC.Range("M80").Value = -200001922
$ws_ALC.Range("N80").Value = -5596

$ws_ALC.Range("H83").Value = 60607056
$ws_ALC.Range("I83").Value = 66667640
$ws_ALC.Range("J83").Value = 1200
$ws_ALC.Range("K83").Value = 600008760
$ws_ALC.Range("L83").Value = 10800
$ws_ALC.Range("M83").Value = -600003768
$ws_ALC.Range("N83").Value = -20784

$ws_ALC.Range("H125").Value = 1266.6666
$ws_ALC.Range("J125").Value = 2178.6667
$ws_ALC.Range("L125").Value = 19608.0003
$ws_ALC.Range("N125").Value = -24528.0003

$ws_ALC.Range("H129").Value = 365166.72
$ws_ALC.Range("I129").Value = 924258.25
$ws_ALC.Range("J129").Value = 1757.25
$ws_ALC.Range("K129").Value = 2772774.75
$ws_ALC.Range("L129").Value = 5271.75
$ws_ALC.Range("M129").Value = -2767774.75
$ws_ALC.Range("N129").Value = -15271.75

$ws_ALC.Range("H137").Value = 2879.9412
$ws_ALC.Range("I137").Value = 989.76666
$ws_ALC.Range("J137").Value = 3910.9456
$ws_ALC.Range("K137").Value = 2969.29998
$ws_ALC.Range("L137").Value = 11732.8368
$ws_ALC.Range("M137").Value = -419.2999799999998
$ws_ALC.Range("N137").Value = -16832.8368

$ws_ALC.Range("H138").Value = 3648.1428
$ws_ALC.Range("I138").Value = 2909.9524
$ws_ALC.Range("J138").Value = 3924.9644
$ws_ALC.Range("K138").Value = 8729.8572
$ws_ALC.Range("L138").Value = 11774.8932
$ws_ALC.Range("M138").Value = -3589.8572
$ws_ALC.Range("N138").Value = -22054.8932

$ws_ALC.Range("H141").Value = 2527.6562
$ws_ALC.Range("I141").Value = 2322.5
$ws_ALC.Range("J141").Value = 5605
$ws_ALC.Range("K141").Value = 6967.5
$ws_ALC.Range("L141").Value = 16815
$ws_ALC.Range("M141").Value = -1787.5
$ws_ALC.Range("N141").Value = -27175

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 43651.2
$ws_ARM.Range("I32").Value = 48501.715
$ws_ARM.Range("K32").Value = 48501.715
$ws_ARM.Range("M32").Value = -48214.715

$ws_ARM.Range("H123").Value = 0
$ws_ARM.Range("J123").Value = 0
$ws_ARM.Range("L123").Value = 0
$ws_ARM.Range("N123").ClearContents()

$ws_ARM.Range("H128").Value = 0
$ws_ARM.Range("J128").Value = 0
$ws_ARM.Range("L128").Value = 0
$ws_ARM.Range("N128").ClearContents()

$ws_ARM.Range("H129").Value = 49999
$ws_ARM.Range("J129").Value = 49999
$ws_ARM.Range("L129").Value = 49999
$ws_ARM.Range("N129").Value = -59999

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H16").Value = 2070.3
$ws_CRP.Range("I16").Value = 1878
$ws_CRP.Range("J16").Value = 2262.6
$ws_CRP.Range("K16").Value = 1878
$ws_CRP.Range("L16").Value = 2262.6
$ws_CRP.Range("M16").Value = -1591
$ws_CRP.Range("N16").Value = -2836.6

$ws_CRP.Range("H31").Value = 19219
$ws_CRP.Range("I31").Value = 0
$ws_CRP.Range("J31").Value = 19219
$ws_CRP.Range("K31").Value = 0
$ws_CRP.Range("L31").Value = 19219
$ws_CRP.Range("M31").ClearContents()
$ws_CRP.Range("N31").Value = -19809

$ws_CRP.Range("H34").Value = 19219
$ws_CRP.Range("I34").Value = 0
$ws_CRP.Range("J34").Value = 19219
$ws_CRP.Range("K34").Value = 0
$ws_CRP.Range("L34").Value = 19219
$ws_CRP.Range("M34").ClearContents()
$ws_CRP.Range("N34").Value = -19623

$ws_CRP.Range("H113").Value = 2070.3
$ws_CRP.Range("I113").Value = 1878
$ws_CRP.Range("J113").Value = 2262.6
$ws_CRP.Range("K113").Value = 1878
$ws_CRP.Range("L113").Value = 2262.6
$ws_CRP.Range("M113").Value = 292
$ws_CRP.Range("N113").Value = -6602.6

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H4").Value = 468.8
$ws_CUL.Range("I4").Value = 468.8
$ws_CUL.Range("J4").Value = 0
$ws_CUL.Range("K4").Value = 1406.4
$ws_CUL.Range("L4").Value = 0
$ws_CUL.Range("M4").Value = -1294.4
$ws_CUL.Range("N4").ClearContents()

$ws_CUL.Range("H107").Value = 11892.056
$ws_CUL.Range("I107").Value = 20718.6
$ws_CUL.Range("J107").Value = 8497.23
$ws_CUL.Range("K107").Value = 62155.8
$ws_CUL.Range("L107").Value = 25491.69
$ws_CUL.Range("M107").Value = -60235.8
$ws_CUL.Range("N107").Value = -29331.69

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H34").Value = 20000
$ws_GSM.Range("J34").Value = 20000
$ws_GSM.Range("L34").Value = 20000
$ws_GSM.Range("N34").Value = -20536

$ws_GSM.Range("H76").Value = 20000
$ws_GSM.Range("J76").Value = 20000
$ws_GSM.Range("L76").Value = 20000
$ws_GSM.Range("N76").Value = -20630

$ws_GSM.Range("H79").Value = 20000
$ws_GSM.Range("J79").Value = 20000
$ws_GSM.Range("L79").Value = 20000
$ws_GSM.Range("N79").Value = -22184

$ws_GSM.Range("H123").Value = 34873.5
$ws_GSM.Range("J123").Value = 34873.5
$ws_GSM.Range("L123").Value = 34873.5
$ws_GSM.Range("N123").Value = -39773.5

$ws_GSM.Range("H126").Value = 11780.818
$ws_GSM.Range("I126").Value = 12788.9
$ws_GSM.Range("K126").Value = 38366.7
$ws_GSM.Range("M126").Value = -35896.7

$ws_GSM.Range("H132").Value = 7571.4287
$ws_GSM.Range("I132").Value = 11168.167
$ws_GSM.Range("J132").Value = 4873.875
$ws_GSM.Range("K132").Value = 33504.501
$ws_GSM.Range("L132").Value = 14621.625
$ws_GSM.Range("M132").Value = -30974.501
$ws_GSM.Range("N132").Value = -19681.625

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 55557924
$ws_LTW.Range("I7").Value = 71430650
$ws_LTW.Range("K7").Value = 71430650
$ws_LTW.Range("M7").Value = -71430538

$ws_LTW.Range("H16").Value = 2660.2222
$ws_LTW.Range("I16").Value = 1705.1666
$ws_LTW.Range("K16").Value = 1705.1666
$ws_LTW.Range("M16").Value = -1535.1666

$ws_LTW.Range("H40").Value = 2551.5625
$ws_LTW.Range("I40").Value = 2511.4546
$ws_LTW.Range("K40").Value = 2511.4546
$ws_LTW.Range("M40").Value = -2375.4546

$ws_LTW.Range("H61").Value = 2620.5417
$ws_LTW.Range("I61").Value = 2690.5908
$ws_LTW.Range("J61").Value = 1850
$ws_LTW.Range("K61").Value = 2690.5908
$ws_LTW.Range("L61").Value = 1850
$ws_LTW.Range("M61").Value = -2488.5908
$ws_LTW.Range("N61").Value = -2254

$ws_LTW.Range("H113").Value = 2620.5417
$ws_LTW.Range("I113").Value = 2690.5908
$ws_LTW.Range("J113").Value = 1850
$ws_LTW.Range("K113").Value = 2690.5908
$ws_LTW.Range("L113").Value = 1850
$ws_LTW.Range("M113").Value = -520.5907999999999
$ws_LTW.Range("N113").Value = -6190

$ws_LTW.Range("H126").Value = 55557924
$ws_LTW.Range("I126").Value = 71430650
$ws_LTW.Range("K126").Value = 214291950
$ws_LTW.Range("M126").Value = -214289480

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H2").Value = 2000
$ws_WVR.Range("J2").Value = 2000
$ws_WVR.Range("L2").Value = 2000
$ws_WVR.Range("N2").Value = -2224

$ws_WVR.Range("H141").Value = 0
$ws_WVR.Range("I141").Value = 0
$ws_WVR.Range("J141").Value = 0
$ws_WVR.Range("K141").Value = 0
$ws_WVR.Range("L141").Value = 0
$ws_WVR.Range("M141").ClearContents()
$ws_WVR.Range("N141").ClearContents()
